$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reviews")

# Row 2 keeps its position; only the score (E) changes.
$ws.Cells.Item(2, 5).Value = 452.753

# Rows 3-16 are replaced in full with the re-ranked/re-ordered data.
# Each line below is one destination row (sheet rows 3..16, in order),
# with 10 tab-separated fields: A,B,C,D,E,F,G,H,I,J.
# Text fields (A,B,C,D,F,G) are base64-encoded (UTF-8) to avoid quoting issues;
# numeric fields (E,H,I,J) are plain decimal numbers.
$rowData = @"
REwtNzM1ZjZmNWU3ZWU2ZDk5OGFlY2RhNTlhYmFiY2RiY2M=	SG93IEkgcGxhbiB0byByZXZpZXcgdGhlIE5pbnRlbmRvIFN3aXRjaCAy	aHR0cHM6Ly93d3cudGhldmVyZ2UuY29tL25pbnRlbmRvLzY3OTM0Ni9uaW50ZW5kby1zd2l0Y2gtMi1yZXZpZXctZ3VpZGVsaW5lcy1wbGFuLWxhdW5jaA==	TmludGVuZG8gaGFzIGRlY2lkZWQgdG8gbm90IHNlbmQgb3V0IGVhcmx5IHJldmlldyB1bml0cyB0byBUaGUgVmVyZ2Ugb3Igb3RoZXIgb3V0bGV0cywgY2l0aW5nIHRoZSBuZWVkIGZvciBkYXktb25lIHNvZnR3YXJlIHVwZGF0ZXMuIEl04oCZcyBhIHdheSB0byBleHBsb3JlIGV2ZXJ5IGZhY2V0IG9mIGEgbmV3IGNvbnNvbGUgYW5kIGdldCBhIHdpZGUgdmFyaWV0eSBvZiBwZXJzcGVjdGl2ZXMuIFRoaXMgdGltZSwgd2XigJlyZSBnb2luZyB0byBzdGFydCB3aXRoIHRob3NlIGRlZXBlciBkaXZlcy4=	253.28822	MjAyNS0wNi0wNFQxMjowMDowMCswMDowMA==	WydBbG1vc3QgZXZlcnl0aGluZyBhYm91dCB0aGUgU3dpdGNoIDIgYWhlYWQgb2YgbGF1bmNoIGhhcyBiZWVuIGEgbGl0dGxlIHdlaXJkLCBmcm9tIHRoZSBjb25mdXNpbmcgbWVzc2FnaW5nIGFuZCBoaWdoIHByaWNpbmcgdG8gdGhlIHVuZm9ydHVuYXRlbHkgdGltZWQgY29ubmVjdGlvbiB3aXRoIHNwaWtpbmcgdGFyaWZmcy4nLCAnSW50ZW50IHRvIGtlZXAgdGhhdCBzdHJhbmdlIGVuZXJneSBnb2luZywgTmludGVuZG8gaGFzIGRlY2lkZWQgdG8gbm90IHNlbmQgb3V0IGVhcmx5IHJldmlldyB1bml0cyB0byBUaGUgVmVyZ2Ugb3Igb3RoZXIgb3V0bGV0cywgY2l0aW5nIHRoZSBuZWVkIGZvciBkYXktb25lIHNvZnR3YXJlIHVwZGF0ZXMuIChUaGlzIHdhcyBub3QgdGhlIGNhc2Ugd2l0aCB0aGUgb3JpZ2luYWwgU3dpdGNoLiknXQ==	0.0	0.0	0.0
REwtNTZiMzQ1NmZhNGQ0YzJjMDMxNmY4Yzc0OGVjMzg5NDI=	TmludGVuZG8gU3dpdGNoIDIgV2VsY29tZSBUb3VyIFJldmlldw==	aHR0cHM6Ly93d3cucm9ja2V0Y2hhaW5zYXcuY29tLmF1L3Jldmlldy9uaW50ZW5kby1zd2l0Y2gtMi13ZWxjb21lLXRvdXItcmV2aWV3Lw==	VGhlIFdlbGNvbWUgVG91ciBpcyBwcmVzZW50ZWQgYXMgYSB2aXJ0dWFsIG11c2V1bSwgYW4gZW5vcm1vdXMgdmlydHVhbCByZXByZXNlbnRhdGlvbiBvZiB0aGUgU3dpdGNoIDIuIFRoZSBjYW1lcmEgaXMgc2V0IGhpZ2ggZnJvbSBhbiBpc29tZXRyaWMgcGVyc3BlY3RpdmUsIHNvbWV0aW1lcyB6b29taW5nIG91dCB0byBlbmNvbXBhc3MgYW4gZW50aXJlIGFyZWEgd2l0aCBwZW9wbGUgYnJvd3NpbmcgdGhlIGV4aGliaXRzIGxpa2UgYW50cy4gVGhlcmXigJlzIGV2ZW4gc29tZSBtaWxkbHkgZnVubnkgZGlhbG9ndWUgYW5kIGxpdHRsZSBzZXQgcGllY2VzIHRocm93biBpbiBhbW9uZyB0aGUgYXR0ZW5kZWVzIHlvdSBjYW4gY2hhdCB0by4=	247.76453	MjAyNS0wNi0wOVQwNDoyMDo0NCswMDowMA==	W10=	0.8126	0.132	0.0
REwtMmE0MThlNmY3NWY2MDJiY2Q1MjEzNDhjYTRlZGMzM2E=	TmludGVuZG8gU3dpdGNoIDIgQ2FtZXJhIHJldmlldyAtIGdvb2QgZm9yIEdhbWVDaGF0IGJ1dCBub3QgbXVjaCBlbHNl	aHR0cHM6Ly93d3cucG9ja2V0dGFjdGljcy5jb20vbmludGVuZG8tc3dpdGNoLTItY2FtZXJhLXJldmlldw==	SSdtIHRhbGtpbmcgYWJvdXQgdGhlIE5pbnRlbmRvIFN3aXRjaCAyIENhbWVyYS4gSXQgZG9lc24ndCBoYXZlIG1hbnkgdXNlcyBvdXRzaWRlIG9mIEdhbWVDaGF0LCB0aGUgU3dpdGNoIDIncyBuZXcgc29jaWFsIGZlYXR1cmUsIGFuZCBhIGZldyBNYXJpbyBQYXJ0eSBKYW1ib3JlZSDigJMgU3dpdGNoIDIgRWRpdGlvbiBtaW5pZ2FtZXMuIEF0IFBvY2tldCBUYWN0aWNzLCBvdXIgZXhwZXJ0cyBzcGVuZCBkYXlzIHRlc3RpbmcgZ2FtZXMsIHBob25lcywgdGVjaCwgYW5kIHNlcnZpY2VzLiBXZSBhbHdheXMgc2hhcmUgaG9uZXN0IG9waW5pb25zIHRvIGhlbHAgeW91IGJ1eSB0aGUgYmVzdC4=	225.1066	MjAyNS0wNi0wN1QxNjo1NTo0NQ==	WyJVbmxpa2UgbW9zdCByZXZpZXdzIG9uIFBvY2tldCBUYWN0aWNzLCBJJ20gbm90IHNjb3JpbmcgdGhlIE5pbnRlbmRvIFN3aXRjaCAyIENhbWVyYSBvdXQgb2YgMTAgcmlnaHQgbm93LiBUaGF0J3MgYmVjYXVzZSBJJ3ZlIG5vdCB1c2VkIGFueSBvdGhlciBTd2l0Y2ggMiBjYW1lcmFzLCBzbyB0aGVyZSdzIG5vIHByaW9yIGNvbnRleHQgdG8gYmFzZSBhIHNjb3JlIG9uLiIsICJGb3IgbW9yZSBvZiBvdXIgcG9zdC1sYXVuY2ggY292ZXJhZ2UsIGJlIHN1cmUgdG8gdGFrZSBhIGxvb2sgYXQgb3VyIE5pbnRlbmRvIFN3aXRjaCAyIHJldmlldywgTmludGVuZG8gU3dpdGNoIDIgUHJvIENvbnRyb2xsZXIgcmV2aWV3LCBhbmQgTWFyaW8gS2FydCBXb3JsZCByZXZpZXcgd2hpbGUgeW91J3JlIGhlcmUuIiwgIk9yLCBpZiB5b3UnZCByYXRoZXIgZ3JhYiBzb21ldGhpbmcgdG8gY29tcGxldGUgeW91ciBuZXcgZ2FtaW5nIHNldHVwLCBzZWUgb3VyIGd1aWRlcyB0byB0aGUgYmVzdCBOaW50ZW5kbyBTd2l0Y2ggMiBhY2Nlc3NvcmllcyBhbmQgdGhlIGJlc3QgTmludGVuZG8gU3dpdGNoIDIgY29udHJvbGxlcnMuIl0=	0.9337	0.247	0.0
REwtNGI4NmU5NjhhOWIxNjQ2ZWY0NzE2YjhiN2MyNmQ5YWQ=	U3RhcnRpbmcgc2hvdCBmb3IgTmludGVuZG8gU3dpdGNoIDIg4oCTIEdlcm1hbiBzdHJlYW1lciBpcyBhbHJlYWR5IHBsYXlpbmcgTWFyaW8gS2FydCBXb3JsZA==	aHR0cHM6Ly93d3cuZ2lnYS5kZS9nYW1lcy9zdGFydHNjaHVzcy1mdWVyLW5pbnRlbmRvLXN3aXRjaC0yLWRldXRzY2hlLXN0cmVhbWVyaW4tem9ja3Qtc2Nob24tbWFyaW8ta2FydC13b3JsZC0tMDFKV1hONjQ1SlBQSlNDTUY3OVpZUTdTUUg=	TmludGVuZG8gU3dpdGNoIDI6IERheS1PbmUgUGF0Y2ggaXMgaGVyZSBUaGUgZmlyc3QgZ2FtZXJzIGFscmVhZHkgaGF2ZSB0aGVpciBTd2l0Y2ggMiBhdCBob21lLiBWZXJzaW9uIDIwLjEuMSBpcyBsaXZlIGFuZCB0aGVyZSBhcmUgZmlyc3QgZXhwZXJpZW5jZXMgd2l0aCB0aGUgbmV3IGh5YnJpZCBjb25zb2xlLiBOaW50ZW5kbyBkaWRuJ3QgbWlzcyBhIGRheSwgaXQncyBhbHJlYWR5IEp1bmUgNXRoIGluIE5ldyBaZWFsYW5kLg==	47.643147	MjAyNS0wNi0wNFQxNDoyMzowOCswMDowMA==	WydOaW50ZW5kbyBTd2l0Y2ggMjogRGF5LU9uZSBQYXRjaCBpcyBoZXJlIFRoZSBmaXJzdCBnYW1lcnMgYWxyZWFkeSBoYXZlIHRoZWlyIFN3aXRjaCAyIGF0IGhvbWUuIFVudGlsIG5vdywgdGhlIGNvbnNvbGUgd2FzIGNvbXBsZXRlbHkgdXNlbGVzcyDigJMgYSBkYXktb25lIHBhdGNoIGZyb20gTmludGVuZG8gd2FzIG1pc3NpbmcuIFRoYXQgaGFzIG5vdyBjaGFuZ2VkLicsICdUaGUgR2VybWFuIHN0cmVhbWVyIExvc3RLaXR0biBoYXMgbm90IG9ubHkgZ290IGhlciBTd2l0Y2ggMiwgc2hlIGlzIGV2ZW4gbm93IHBsYXlpbmcgTWFyaW8gS2FydCBXb3JsZCBvbiBUd2l0Y2guIEFjY29yZGluZyB0byB0aGUgcGlubmVkIG1lc3NhZ2Ugb2YgYSBjaGF0IG1vZGVyYXRvciwgc2hlIGhhcyByZWNlaXZlZCB0aGUgY29uc29sZSBmcm9tIE5pbnRlbmRvIGFuZCBpcyBhbGxvd2VkIHRvIHBsYXkgc2luY2UgMiBwbS4nLCAnUGxheWVycyBzaG93Y2FzZSBTd2l0Y2ggMiBlU2hvcCBPdGhlciBwbGF5ZXJzIGFsc28gc2hhcmUgdGhlaXIgZXhwZXJpZW5jZXMgd2l0aCBTd2l0Y2ggMi4gT24gUmVkZGl0LCB1c2VycyBzaG93IEdhYlNhbSwgZm9yIGV4YW1wbGUsIHRoZSBuZXcgYW5kIGltcHJvdmVkIE5pbnRlbmRvIGVTaG9wLiBJdCBub3cgcnVucyBtdWNoIG1vcmUgc21vb3RobHkgdGhhbiB0aGUgb25lIG9uIHRoZSBvbGQgY29uc29sZS4nXQ==	-0.1531	0.0	0.039
REwtYjQ4MzVmNGZlMWE5NmVhMmQ5ZjNlNWE4MzQ0OWNlZmM=	SSd2ZSBhbHJlYWR5IHRyaWVkIHRoZSdtb3VzZScgbW9kZSBvbiBOaW50ZW5kbyBTd2l0Y2ggMiwgaXQgb3BlbnMgYSBkb29yIHRoYXQgaGFzIHNvIGZhciBiZWVuIGNsb3NlZCBvbiBjb25zb2xlcy4=	aHR0cHM6Ly93d3cueGF0YWthLmNvbS92aWRlb2p1ZWdvcy9oZS1wcm9iYWRvLW1vZG8tcmF0b24tbmludGVuZG8tc3dpdGNoLTItYWJyZS1wdWVydGEtYWhvcmEtY2VycmFkYS1jb25zb2xhcw==	TmludGVuZG8gU3dpdGNoIDIgaXMgbm90IGEgbWFjaGluZSB0byBleHBlcmltZW50IHdpdGguIFRoZSAxNTIgbWlsbGlvbiBTd2l0Y2hlcyBzb2xkIGFyZSBhIHJlZmxlY3Rpb24gdGhhdCB0aGUgY29tcGFueSBoYXMgYWNoaWV2ZWQgYSBwcm9kdWN0IHRoYXQgcGVvcGxlIGhhdmUgbGlrZWQuIFRoZXkgZG9uJ3QgbmVlZCBzb21ldGhpbmcgYXMgZGlzcnVwdGl2ZSBhcyBTd2l0Y2ggb3IgV2lpIGluIHRoZWlyIGRheSwgYnV0IGEgY29udGludW91cyBjb25zb2xlLiBBbmQgYXBhcnQgZnJvbSB0aGUgc2l6ZSwgYSBjaGFuZ2UgY29tZXMgd2l0aCB0aGUgb3B0aWNhbCB0ZWNobm9sb2d5Lg==	41.942944	MjAyNS0wNi0wNVQxNjoxNTo0OQ==	WydUdXJuaW5nIHRvIHRoZSBKb3ktQ29uIE5pbnRlbmRvIFN3aXRjaCAyIGlzIG5vdCBhIG1hY2hpbmUgdG8gZXhwZXJpbWVudCB3aXRoLiBUaGUgMTUyIG1pbGxpb24gU3dpdGNoZXMgc29sZCwgYW5kIHJpc2luZywgYXJlIGEgcmVmbGVjdGlvbiB0aGF0IHRoZSBjb21wYW55IGhhcyBhY2hpZXZlZCBhIHByb2R1Y3QgdGhhdCBwZW9wbGUgaGF2ZSBsaWtlZC4nXQ==	0.0644	0.036	0.031
REwtZmEyNjllOGEwZTE5YTMyYjJiZjBmOGNlNzMxZDgyYjE=	TWFyaW8gS2FydCBXb3JsZCBvbiBNZXRhY3JpdGljOiBGaXJzdCB0ZXN0IHNjb3JlcyBhcmUgbGl2ZSBhbmQgYW4gYWJzb2x1dGUgZHJlYW0h	aHR0cHM6Ly93d3cuZ2FtZXByby5kZS9hcnRpa2VsL21hcmlvLWthcnQtd29ybGQtbWV0YWNyaXRpYy1lcnN0ZS10ZXN0cy13ZXJ0dW5nLDM0MzQyNDYuaHRtbA==	TmludGVuZG8gU3dpdGNoIDIgaXMgaGVyZSBhbmQgd2UncmUgc2hvb3Rpbmcgb3VyIGZpcnN0IDM3IG1pbnV0ZXMgaW4gTWFyaW8gS2FydCBXb3JsZCBBdXRvcGxheSBXaGF0J3MgYmVpbmcgcHJhaXNlZD8gTmludGVuZHVvIGdpdmVzIGl0IGEgc2NvcmUgb2YgOTIgYW5kIHdyaXRlczogIkV2ZXJ5dGhpbmcgY2FuIHN0aWxsIGNoYW5nZSIgTmludGVuZG8gaGFzIHBlcmZlY3RseSBiYWxhbmNlZCBjbGFzc2ljIGFuZCBtb2Rlcm4gZWxlbWVudHMgdG8gY3JlYXRlIGFuIGVuZ2FnaW5nIGFuZCBhY2Nlc3NpYmxlIGV4cGVyaWVuY2UgdGhhdCB3aWxsIGRlbGlnaHQgYm90aCB2ZXRlcmFucyBhbmQgbmV3Y29tZXJzLg==	39.233864	MjAyNS0wNi0wNlQxMDoxNDowMA==	WyJUaGUgZmlyc3QgcmF0aW5ncyBmb3IgTWFyaW8gS2FydCBXb3JsZCBhcmUgZW50aHVzaWFzdGljLiBObyBOaW50ZW5kbyBTd2l0Y2ggMiBjb25zb2xlcyBhbmQgcmV2aWV3IHNhbXBsZXMgd2VyZSBzZW50IGluIGFkdmFuY2UsIHdlIHJlY2VpdmVkIE5pbnRlbmRvJ3MgbmV3ZXN0IGZsYWdzaGlwIG9uIHJlbGVhc2UgZGF5LiIsICIzNzo0OCBUaGUgTmludGVuZG8gU3dpdGNoIDIgaXMgaGVyZSBhbmQgd2UncmUgc2hvb3Rpbmcgb3VyIGZpcnN0IDM3IG1pbnV0ZXMgaW4gTWFyaW8gS2FydCBXb3JsZCBBdXRvcGxheSBXaGF0J3MgYmVpbmcgcHJhaXNlZD8iLCAnVGhpcyBtZWFuczogTWFyaW8gS2FydCBXb3JsZCBpcyB0aGUgbGF0ZXN0IGFuZCBwZXJoYXBzIGJlc3QgcGFydCBvZiB0aGUgbG9uZy1zdGFuZGluZyBhcmNhZGUgcmFjaW5nIHNlcmllcyBhbmQgaXMgdGFrZW4gdG8gdGhlIG5leHQgbGV2ZWwgd2l0aCB0aGUgaGFyZHdhcmUgb2YgTmludGVuZG8gU3dpdGNoIDIgYW5kIGFuIGV4dGVuc2l2ZSByb3N0ZXIgb2YgZHJpdmVycywgdHJhY2tzIGFuZCBlbmRsZXNzIGZ1bi4nXQ==	0.9413	0.252	0.0
REwtMjJhZTY2ZmMxZjA5OWYwMjBkNjk4ODZiMTQ0MWU1NmU=	TmludGVuZG8gc2VsbHMgMy41IG1pbGxpb24gU3dpdGNoIDIgY29uc29sZXMgaW4gZm91ciBkYXlz	aHR0cHM6Ly9nYW1lbGluZXIubmwvbmlldXdzL25pbnRlbmRvLXZlcmtvb3B0LTM1LW1pbGpvZW4tc3dpdGNoLTItY29uc29sZXMtYmlubmVuLXZpZXItZGFnZW4vNTEwMjA=	TmludGVuZG8gU3dpdGNoIDIgaXMgdGhlIG1vc3Qgc3VjY2Vzc2Z1bCBsYXVuY2ggaW4gdGhlIGhpc3Rvcnkgb2YgTmludGVuZG8gaGFyZHdhcmUuIFRoZSBkZXNpcmVkIGNvbnNvbGUgc29sZCBvdmVyIDMuNSBtaWxsaW9uIGNvcGllcyBpbiBmb3VyIGRheXMuIFRoaXMgaXMgYW4gaW5kaWNhdGlvbiBvZiB0aGUgdHJ1c3QgTmludGVuZG8gaGFzIGluIHByb2R1Y2luZyBhIGNvbnNvbGUgdGhhdCBpcyBtb3JlIHRoYW4gd29ydGggdGhlIGVmZm9ydC4=	37.532463	MjAyNS0wNi0xMVQwODoxNTowMCswMDowMA==	WydBbmQgb2YgY291cnNlLCB0aGlzIGlzIGFsc28gYW4gaW5kaWNhdGlvbiBvZiB0aGUgdHJ1c3QgdGhhdCBOaW50ZW5kbyBoYXMgaW4gcHJvZHVjaW5nIGEgY29uc29sZSB0aGF0IGlzIG1vcmUgdGhhbiB3b3J0aCB0aGUgZWZmb3J0LicsICJUaGUgMy41IG1pbGxpb24gdW5pdHMgc29sZCBpbiBmb3VyIGRheXMgbWFrZSB0aGUgTmludGVuZG8gU3dpdGNoIDIgdGhlIG1vc3Qgc3VjY2Vzc2Z1bCBsYXVuY2ggaW4gdGhlIGhpc3Rvcnkgb2YgTmludGVuZG8gaGFyZHdhcmUuIENoZWNrIG91dCBvdXIgTmludGVuZG8gU3dpdGNoIDIgcmV2aWV3IGhlcmUuIEluIGFkZGl0aW9uLCB3ZSd2ZSBnb3QgYSBNYXJpbyBLYXJ0IFdvcmxkIHJldmlldyByZWFkeSBmb3IgeW91LiJd	0.8927	0.235	0.0
REwtYzU2MWI0ZWIzNzY0MTk1YTFmZDRlMmVhMWI5ZTU4MWI=	RnJlZSB1cGRhdGUgdG8gUG9rw6ltb24gU2NhcmxldCBhbmQgVmlvbGV0IGZvciBOaW50ZW5kbyBTd2l0Y2ggMg==	aHR0cHM6Ly90aGF0c2dhbWluZy5ubC9ncmF0aXMtdXBkYXRlLXBva2Vtb24tc2NhcmxldC1lbi12aW9sZXQtdm9vci1uaW50ZW5kby1zd2l0Y2gtMi8=	UG9rw6ltb24gU2NhcmxldCBhbmQgUG9rw6ltb24gVmlvbGV0IGFyZSBhdmFpbGFibGUgZm9yIHRoZSBOaW50ZW5kbyBTd2l0Y2ggMi4gVGhpcyB1cGRhdGUgbWFrZXMgdGhlIGFkdmVudHVyZXMgaW4gdGhlIFBhbGRlYSByZWdpb24gZXZlbiBtb3JlIGV4Y2l0aW5nLiBJbXByb3ZlZCBpbWFnZSBxdWFsaXR5IGFsbG93cyBmYW5zIHRvIHJlY2hhcmdlIHRoZWlyIGJhdHRlcmllcy4=	34.85138	MjAyNS0wNi0wNVQwOTozODo1MiswMDowMA==	WyJBIHNoYXJwZXIgaW1hZ2UgcXVhbGl0eTogdGhlIG9wdGltaXphdGlvbnMgZm9yIGJvdGggdGhlIE5pbnRlbmRvIFN3aXRjaCAyJ3Mgc2NyZWVuIGFuZCBIRCB0ZWxldmlzaW9ucyBwcm92aWRlIGV2ZW4gYmV0dGVyIGltYWdlIHF1YWxpdHkuIEEgaGlnaGVyIGZyYW1lIHJhdGU6IHRoZSBpbXByb3ZlZCBmcmFtZSByYXRlcyBvZiB0aGUgTmludGVuZG8gU3dpdGNoIDIgYWxsb3cgZm9yIG1vcmUgZmx1aWQgbW90aW9uLiIsICdBbiBpbW1lcnNpdmUgUGFsZGVhOiBvbiB0aGUgTmludGVuZG8gU3dpdGNoIDIsIHRoZSBhbHJlYWR5IGJ1c3RsaW5nIFBhbGRlYSByZWdpb24gY29tZXMgdG8gbGlmZSBpbiBhbiB1bnByZWNlZGVudGVkIHdheSwgd2hldGhlciBwbGF5ZXJzIGV4cGxvcmUgdGhlIG9wZW4gd29ybGQsIGZpZ2h0IFBva8OpbW9uIG9yIGNvbXBsZXRlIHRoZSBQb2vDqWRleC4nLCAiVHJhaW5lcnMgd2hvIGhhdmUgYSBOaW50ZW5kbyBTd2l0Y2ggMiBhbmQgUG9rw6ltb24gU2NhcmxldCBvciBQb2vDqW1vbiBWaW9sZXQgY2FuIG5vdyBleHBlcmllbmNlIHRoZSBnYW1lcyBldmVuIG1vcmUgaW50ZW5zZWx5IGFuZCBiZWF1dGlmdWxseS4gQW5kIGl0J3MgY29tcGxldGVseSBmcmVlLiAwIGNvbW1lbnRzIDAiXQ==	0.8508	0.23	0.0
REwtMmZhMjZiZjcyNDg1MmFiYzVjNzJjOWY4MjIwMjhjNTI=	WmVsZGE6IFRoZSB1bmRlcnJhdGVkIGdhbWUgdGhhdCBkZXNlcnZlcyBhIHJlbWFrZSBmb3IgU3dpdGNoIDI=	aHR0cHM6Ly93d3cuemF6b29tLml0LzIwMjUtMDYtMDcvemVsZGEtaWwtZ2lvY28tc290dG92YWx1dGF0by1jaGUtbWVyaXRhLXVuLXJlbWFrZS1wZXItc3dpdGNoLTIvMTcxNDQ3ODEv	TmludGVuZG8gd291bGQgaGF2ZSB0aGUgb3Bwb3J0dW5pdHkgdG8gc3RyZW5ndGhlbiB0aGUgY29ubmVjdGlvbiBiZXR3ZWVuIHRoZSBwYXN0IGFuZCB0aGUgZnV0dXJlLCBwcm92aWRpbmcgYW4gdW5mb3JnZXR0YWJsZSBleHBlcmllbmNlLiBXaXRoIHRoZSBhZHZlbnQgb2YgdGhlIG5ldyBjb25zb2xlLCB0aGUgY2hhbGxlbmdlIGZvciBOaW50ZW5kbyB3aWxsIGJlIHRvIGNvbnRpbnVlIHRvIGlubm92YXRlIGFuZCBkZWxpdmVyIGV4cGVyaWVuY2VzIHRoYXQgbWVldCB0aGUgZXhwZWN0YXRpb25zIG9mIHRoZSBwYXNzaW9uYXRlLiBUaGUgTGVnZW5kIG9mIFplbGRhOiBUZWFycyBvZiB0aGUgS2luZ2RvbSBnYW1lIHNob3dzIHVwIG9uIE5pbnRlbmRvIFN3aXRjaCAyIHdpdGggMjAgbWludXRlcyBvZiBnYW1lcGxheS4=	24.060257	MjAyNS0wNi0wN1QxODoxNjoyOA==	W10=	0.8779	0.178	0.026
REwtMDU2NTBkZmI3OWJkOGJmN2VmOTk3YmIyY2UyZjNjZDM=	TmludGVuZG8gU3dpdGNoIDIgRmlyc3QgSW1wcmVzc2lvbnMgQW5kIEEgU3VycHJpc2luZyBQdXJjaGFzaW5nIFNjb3Jl	aHR0cHM6Ly9ob3RoYXJkd2FyZS5jb20vbmV3cy9uaW50ZW5kby1zd2l0Y2gtMi1maXJzdC1pbXByZXNzaW9ucy1hbmQtYS1zdXJwcmlzaW5nLXB1cmNoYXNpbmc=	TmludGVuZG8gU3dpdGNoIDIgaXMgZmluYWxseSBvdXQsIGFuZCBJIGhhdmUgc29tZSBmaXJzdCBpbXByZXNzaW9ucyBvbiB0aGUgZXhwZXJpZW5jZSBvZiBib3RoIHBsYXlpbmcgYW5kIHB1cmNoYXNpbmcgdGhlIGNvbnNvbGUuIEkgd2FzIG5vdCBsdWNreSBlbm91Z2ggdG8gc2VjdXJlIGEgcHJlb3JkZXIgYWhlYWQgb2YgdGhlIGNvbnNvbGUncyBKdW5lIDV0aCByZWxlYXNlLiBUaGlyZC1wYXJ0eSBlbnRyZXByZW5ldXJzIHdlcmUgYWxyZWFkeSBsaXN0aW5nIHRoZSBjb25zb2xlIGZvciBhIHByZW1pdW0gb3ZlciBpdHMgJDQ0OSBNU1JQLiBNeSBsb2NhbCBXYWxtYXJ0IHdhcyBuZWFybHkgZW1wdHkgYW5kIHVwb24gaW5xdWlyaW5nLCBJIHdhcyB0b2xkIGl0IGhhZCBhcm91bmQgc2l4IHRvIGVpZ2h0IGNvbnNvbGVzIHN0aWxsIGF2YWlsYWJsZSBpbiB0aGUgYmFjay4=	21.127535	MjAyNS0wNi0wNlQxNDozMTowMCswMDowMA==	WydUaGUgTmludGVuZG8gU3dpdGNoIDIgaXMgZmluYWxseSBvdXQsIGFuZCBJIGhhdmUgc29tZSBmaXJzdCBpbXByZXNzaW9ucyBvbiB0aGUgZXhwZXJpZW5jZSBvZiBib3RoIHBsYXlpbmcgYW5kIHB1cmNoYXNpbmcgdGhlIGNvbnNvbGUuJ10=	0.5673	0.114	0.054
REwtOGQ4Y2FmMDRkNjVkYThlYmFjNjEzNmI3Y2VjNTlhMTE=	TmludGVuZG8gU3dpdGNoIDIgVW5ib3hpbmcsIFNldHVwICYgV2hhdCBZb3UgTmVlZCB0byBLbm93	aHR0cHM6Ly93d3cuZ2Vla3ktZ2FkZ2V0cy5jb20vbmludGVuZG8tc3dpdGNoLTItNC8=	VGhpcyBndWlkZSBleHBsb3JlcyB0aGUgdW5ib3hpbmcgZXhwZXJpZW5jZSwgc2V0dXAgcHJvY2VzcywgZGVzaWduIGNoYW5nZXMsIGhhcmR3YXJlIGFkdmFuY2VtZW50cywgYW5kIG5ldyBmZWF0dXJlcywgd2hpbGUgYWxzbyBjb21wYXJpbmcgaXQgdG8gdGhlIG9yaWdpbmFsIE5pbnRlbmRvIFN3aXRjaC4gVGhlIHBhY2thZ2luZyBpcyBjb21wYWN0IGFuZCBlY28tZnJpZW5kbHksIHVuZGVyc2NvcmluZyBOaW50ZW5kb+KAmXMgY29tbWl0bWVudCB0byByZWR1Y2luZyBlbnZpcm9ubWVudGFsIGltcGFjdC4gSW5zaWRlIHRoZSBib3gsIHlvdeKAmWxsIGZpbmQ6IFRoZSBOaW50ZW5kbyBTd2l0Y2ggMiBjb25zb2xlIFR3byByZWRlc2lnbmVkIEpveS1Db24gY29udHJvbGxlcnMgQSByZXZhbXBlZCBkb2NrIGZlYXR1cmluZyBhbiBFdGhlcm5ldCBwb3J0IEFuIEhETUkgY2FibGUgQSBwb3dlciBhZGFwdGVyLg==	19.452656	MjAyNS0wNi0wOFQxMzowMDoyNA==	WydUaGUgdmlkZW8gYmVsb3cgZnJvbSBab2xsb3RlY2ggZ2l2ZXMgdXMgYSBkZXRhaWxlZCBsb29rIGF0IHRoZSBuZXcgU3dpdGNoIGNvbnNvbGUuXG5VbmJveGluZzogV2hhdOKAmXMgSW5zaWRlIHRoZSBCb3g/XG5UaGUgdW5ib3hpbmcgZXhwZXJpZW5jZSBvZiB0aGUgTmludGVuZG8gU3dpdGNoIDIgcmVmbGVjdHMgYSBmb2N1cyBvbiBzaW1wbGljaXR5IGFuZCBzdXN0YWluYWJpbGl0eS4nLCAnU2V0dGluZyBVcDogRnJvbSBCb3ggdG8gR2FtZXBsYXlcblNldHRpbmcgdXAgdGhlIE5pbnRlbmRvIFN3aXRjaCAyIGlzIGEgc3RyYWlnaHRmb3J3YXJkIHByb2Nlc3MgZGVzaWduZWQgdG8gZ2V0IHlvdSBnYW1pbmcgcXVpY2tseS4nLCAnTmV3IEZlYXR1cmVzOiBFbmhhbmNpbmcgdGhlIEV4cGVyaWVuY2VcblRoZSBOaW50ZW5kbyBTd2l0Y2ggMiBpbnRyb2R1Y2VzIHNldmVyYWwgbmV3IGZlYXR1cmVzIHRoYXQgZW5oYW5jZSB0aGUgb3ZlcmFsbCB1c2VyIGV4cGVyaWVuY2UgYW5kIGNhdGVyIHRvIG1vZGVybiBnYW1pbmcgbmVlZHM6XG5PbmxpbmUgRnVuY3Rpb25hbGl0eTogRmFzdGVyIGRvd25sb2FkIHNwZWVkcyBhbmQgYSBtb3JlIHN0YWJsZSBjb25uZWN0aW9uIGltcHJvdmUgbXVsdGlwbGF5ZXInXQ==	0.8979	0.166	0.0
REwtZjk1NWMyNjAxMTNjYjQ2MWJjODNkMDdiN2VmNTEzZGQ=	TmludGVuZG8gU3dpdGNoIDIgQ29tZXMgV2l0aCBhIEJlbG93LUF2ZXJhZ2UgRGlzcGxheSBhbmQgRGlzYXBwb2ludGluZyBIRFIgU3VwcG9ydCwgTmV3IEluLURlcHRoIEFuYWx5c2lzIFJldmVhbHM=	aHR0cHM6Ly93Y2NmdGVjaC5jb20vbmludGVuZG8tc3dpdGNoLTItYmVsb3ctYXZlcmFnZS1oZHIv	VGhlIG5ldyBhbmFseXNpcyBieSBHYW1pbmdUZWNoIHRha2VzIGEgZ29vZCBsb29rIGF0IHRoZSBuZXcgTmludGVuZG8gc3lzdGVtJ3MgZGlzcGxheS4gTWVhc3VyaW5nIHRoZSBicmlnaHRuZXNzIGluIFplbGRhOiBCcmVhdGggb2YgdGhlIFdpbGQgcmV0dXJuZWQgYSBtYXhpbXVtIHZhbHVlIG9mIDQyMCBuaXRzLCBidXQgaW4gQ3liZXJwdW5rIDIwNzcsIGl0IGlzIGxvY2tlZCB0byA0NTAgbml0LCBsaWtlbHkgdGhlIG1heGltdW0gcGVhayBicmlnaHRuZXNzIHRoZSBkaXNwbGF5IGlzIGNhcGFibGUgb2YuIFRoaXMgcGVhayBicmlnaHRuZXNzIGlzIG5vdCBldmVuIGNsb3NlIHRvIHByb3ZpZGluZyBhIHByb3BlciBIRFIgZXhwZXJpZW5jZS4=	15.894136	MjAyNS0wNi0wNVQxMTozNzowMCswMDowMA==	WyJSZWxhdGVkIFN0b3J5IExldmVsIFVwIFlvdXIgU3dpdGNoIDIgRXhwZXJpZW5jZSBXaXRoIE11bWJh4oCZcyBCbGFkZSBTZXJpZXMsIENsZWFyIFNlcmllcyBEb2NrYWJsZSBDYXNlcywgYW5kIENhcnJ5aW5nIENhc2VcblRoZSBOaW50ZW5kbyBTd2l0Y2ggMidzIEhEUiBzdXBwb3J0IGluIGRvY2tlZCBtb2RlIGZhcmVzIGEgbGl0dGxlIGJldHRlci4iLCAnVGhlIHN5c3RlbS1sZXZlbCBjYWxpYnJhdGlvbiBpcyBmb3VuZCB0byBiZSBhZGVxdWF0ZSwgYW5kIGluIGdhbWVzIGxpa2UgQ3liZXJwdW5rIDIwNzcsIHRoZSBzeXN0ZW0gZGVsaXZlcnMgYW4gSERSIGV4cGVyaWVuY2Ugb24gcGFyIHdpdGggdGhhdCBvZiB0aGUgb3RoZXIgdmVyc2lvbnMgb2YgdGhlIGdhbWUuJywgJ0luIG90aGVyIGdhbWVzIGxpa2UgWmVsZGE6IEJyZWF0aCBvZiB0aGUgV2lsZCwgdGhlIEhEUiBleHBlcmllbmNlIGlzIHJhdGhlciBkaXNhcHBvaW50aW5nLCBhcyB0aGUgZ2FtZSBsb29rcyB3YXNoZWQgb3V0IGR1ZSB0byBpdHMgYWVzdGhldGljcyBhbmQgdGhlIGxhY2sgb2YgY29udHJhc3QuXG5UaGUgTmludGVuZG8gU3dpdGNoIDIgbGF1bmNoZXMgdG9kYXkgd29ybGR3aWRlLidd	0.928	0.231	0.0
REwtZGFhNzM3MGUxN2JiZjY2OTk5NzdmOTBkZDlkY2MzMmU=	TmludGVuZG8gU3dpdGNoIDI6IE5ldyBnYW1lIGNvbnNvbGUgb2ZmaWNpYWxseSBsYXVuY2hlcyB0b2RheQ==	aHR0cHM6Ly93d3cuYXBwZ2VmYWhyZW4uZGUvbmludGVuZG8tc3dpdGNoLTItbmV1ZS1zcGllbGtvbnNvbGUtZ2VodC1oZXV0ZS1vZmZpemllbGwtYW4tZGVuLXN0YXJ0LTM4MDQwMy5odG1s	TmludGVuZG8gb2ZmaWNpYWxseSB1bnZlaWxlZCB0aGUgbmV3IGdlbmVyYXRpb24gb2YgdGhlIHBvcHVsYXIgU3dpdGNoIGdhbWUgY29uc29sZSwgdGhlIE5pbnRlbmRvIFN3aXRjaCAyLiBUaGUgbmV3IG1vZGVsLCBlcXVpcHBlZCB3aXRoIGEgNy454oCy4oCyIHNjcmVlbiBhbmQgc3VwcG9ydCBmb3IgNEsgY29udGVudCwgaXMgbm93IGF2YWlsYWJsZSBvbiB0aGUgbWFya2V0LiBBdCB0aGUgc2FtZSB0aGlja25lc3MsIHRoZSBjb25zb2xlIGZlYXR1cmVzIGEgc2lnbmlmaWNhbnRseSBsYXJnZXIgNy45LWluY2ggTENEIHNjcmVlbiwgMTA4MHAgcmVzb2x1dGlvbiwgSERSIHN1cHBvcnQsIGFuZCBmcmFtZSByYXRlcyBvZiB1cCB0byAxMjAgZnBzLiBJZiB0aGUgU3dpdGNoIDIgaXMgY29ubmVjdGVkIHRvIGEgVFYgdmlhIHRoZSBhbHNvIHVwZGF0ZWQgZG9jay4gRm9yIDRLIHBsYXliYWNrLCB0aGUgZnJhbWU=	15.065298	MjAyNS0wNi0wNVQxMjozODoxMSswMDowMA==	WydCdXQgd2hhdCBjYW4gdGhlIG5ldyBOaW50ZW5kbyBTd2l0Y2ggMiBkbz8gQ29tcGFyZWQgdG8gaXRzIHByZWRlY2Vzc29yLCB0aGUgTmludGVuZG8gU3dpdGNoIChPTEVEKSwgdGhlIFN3aXRjaCAyIGhhcyBkb25lIHF1aXRlIGEgYml0LicsICdBdCB0aGUgc2FtZSB0aGlja25lc3MsIHRoZSBjb25zb2xlIGZlYXR1cmVzIGEgc2lnbmlmaWNhbnRseSBsYXJnZXIgNy45LWluY2ggTENEIHNjcmVlbiwgMTA4MHAgcmVzb2x1dGlvbiwgSERSIHN1cHBvcnQsIGFuZCBmcmFtZSByYXRlcyBvZiB1cCB0byAxMjAgZnBzLidd	0.802	0.111	0.0
REwtY2NlNzM4N2ZjY2FkNzJhY2Q3YjU3N2RlMDdhMWE2ZTM=	TmludGVuZG8gU3dpdGNoIDIgdW5ib3hpbmc6IFdoYXQgZG8geW91IGdldCB3aXRoIHRoZSBjb25zb2xlPw==	aHR0cHM6Ly93d3cuaW5keTEwMC5jb20vZ2FtaW5nL25pbnRlbmRvLXN3aXRjaC0yLXVuYm94aW5nLWNvbnNvbGUtMjY3MjMxNzQwNQ==	TmludGVuZG8gaGFzIHNlbnQgcmV2aWV3IHVuaXRzIG91dCB0byBtZWRpYSBvdXRsZXRzIGZvciB0aGVtIHRvIGNoZWNrIG91dCB0aGUgbmV3IGNvbnNvbGUgYW5kIGdhbWVzIHJ1bm5pbmcgb24gaXQuIFRoZSBmaXJzdCB0aGluZyB0aGF0IHdpbGwgYmUgc2VlbiB3aGVuIHVuYm94aW5nIGFyZSB0aGUgdHdvIG5ldyBKb3ktQ29uIGNvbnRyb2xsZXJzIGFuZCBtYWluIHBhcnQgb2YgdGhlIGNvbnNvbGUgd2l0aCB0aGUgYmlnZ2VyIHNjcmVlbi4gQmVsb3cgdGhhdCBvbiB0aGUgbmV4dCBsYXllciBkb3duIGlzIGEgYnJpZWYgdXNlciBtYW51YWwsIGEgSERNSSBjYWJsZSB0byBjb25uZWN0IHRoZSBkb2NrIHRvIGEgVFYgYW5kIHRoZSB3cmlzdCBzdHJhcCBhY2Nlc3NvcmllcyB0aGF0IGNvbm5lY3QgdG8gdGhlIFN3aXRjaCAyLg==	14.92691	MjAyNS0wNi0wNVQxNToxOTozNw==	WydOaW50ZW5kbyBoYXMgc2VudCBTd2l0Y2ggMiByZXZpZXcgdW5pdHMgb3V0IHRvIG1lZGlhIG91dGxldHMgZm9yIHRoZW0gdG8gY2hlY2sgb3V0IHRoZSBuZXcgY29uc29sZSBhbmQgZ2FtZXMgcnVubmluZyBvbiBpdC5cbmluZHkxMDAgaXMgdmVyeSBmb3J0dW5hdGUgdG8gaGF2ZSBiZWVuIHNlbnQgb25lIG9mIHRoZXNlIHVuaXRzIGFuZCB0aGlzIGlzIHdoYXQgY29tZXMgaW4gdGhlIGJveCAtIHJlYWQgb3VyIGZpcnN0IGltcHJlc3Npb25zIG9mJ10=	0.5859	0.048	0.0
"@

$lines = $rowData -split "`n"
$destRow = 3
foreach ($line in $lines) {
    $line = $line.Trim("`r")
    if ($line.Length -eq 0) { continue }
    $fields = $line.Split("`t")

    $ws.Cells.Item($destRow, 1).Value = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String($fields[0]))
    $ws.Cells.Item($destRow, 2).Value = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String($fields[1]))
    $ws.Cells.Item($destRow, 3).Value = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String($fields[2]))
    $ws.Cells.Item($destRow, 4).Value = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String($fields[3]))
    $ws.Cells.Item($destRow, 5).Value = [double]$fields[4]
    $ws.Cells.Item($destRow, 6).Value = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String($fields[5]))
    $ws.Cells.Item($destRow, 7).Value = [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String($fields[6]))
    $ws.Cells.Item($destRow, 8).Value = [double]$fields[7]
    $ws.Cells.Item($destRow, 9).Value = [double]$fields[8]
    $ws.Cells.Item($destRow, 10).Value = [double]$fields[9]

    $destRow++
}

# Update Metadata!A2 (count) from 86 to 91.
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("A2").Value = 91
